$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1) Simple text change in the "Developed" section.
Replace-Text `
    "- CRUD update: product/customer/order create+update+delete API wiring done end-to-end." `
    "- Product/Customer/Order create+update+delete API wiring done end-to-end."

# 2) "New in this update" section: 5 bullets collapse into 4 new bullets.
#    Update the text of the first four paragraphs in place, then delete the fifth paragraph entirely.
Replace-Text `
    "- Fixed GitHub Actions CI failure (``MSB1003``) caused by invalid desktop workflow template." `
    "- Added backend store update endpoint: ``PUT /api/stores/{id}`` with owner/admin policy + tenancy checks."

Replace-Text `
    "- Replaced workflow with backend .NET CI:" `
    "- Wired Settings General tab save button to backend store update API."

Replace-Text `
    "  - ``.github/workflows/dotnet-ci.yml``" `
    "- Added save status feedback on Settings page."

Replace-Text `
    "  - Runs restore/build/test in ``backend-dotnet`` project explicitly." `
    "- Backend build verified successful (``dotnet build``) with warnings only."

# Delete the now-redundant fifth paragraph ("Verified locally: ...").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Verified locally*backend-dotnet.csproj*exits successfully.*") {
        $p.Range.Delete()
        break
    }
}

# 3) "Pending / Partial" section: insert a new bullet before the "advanced form fields" bullet,
#    and tweak the wording of that existing bullet.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Some advanced form fields are minimal for now*") {
        $newRange = $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newParaRange = $newPara.Range
        $trimmed = $d.Range($newParaRange.Start, $newParaRange.End - 1)
        $trimmed.Text = "- Store email/phone/address fields are currently UI-only placeholders (not persisted in current store model)."
        break
    }
}

Replace-Text `
    "- Some advanced form fields are minimal for now (e.g., rich item editing in orders, full address management)." `
    "- Some advanced form fields are minimal for now (e.g., rich order item editing, full address management UI)."

# 4) "Git state" section updates.
Replace-Text "- Last pushed commit: 6fe7606" "- Last pushed commit: 131e423"
Replace-Text `
    "- Current CI workflow fix is local and not pushed yet." `
    "- Current store settings API wiring is local and not pushed yet."
